$wb = $excel.ActiveWorkbook

$wsScheme = $wb.Worksheets.Item("Scheme Summary")
$wsUnitSummary = $wb.Worksheets.Item("UnitSummary")
$wsUnitMix = $wb.Worksheets.Item("UnitMix")

# Clear out all the data that drove the UnitSummary and UnitMix sheets -
# the "Scheme Summary" sheet pulls everything via formulas, so clearing the
# source tables collapses those formulas down to "" / #VALUE! results.
$wsUnitSummary.Activate()
$wsUnitSummary.Cells.ClearContents()
$wsUnitSummary.Range("A1:G11").Select()

$wsUnitMix.Activate()
$wsUnitMix.Cells.ClearContents()
$wsUnitMix.Range("D19").Select()

# Leave the original sheet active/selected, matching the recorded UI state.
$wsScheme.Activate()
$wsScheme.Range("A4").Select()
